# [PV-94][WIP] Support for plans without sticky-ids or levels
# Rename header columns on the main plan sheet:
#   "Unique Sticky ID" -> "Row ID"
#   "Task Name"        -> "Task"
#   "Start"             -> "Start Date"
#   "Finish"            -> "End Date"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PV-Test-03-t05-end-date")

$ws.Range("A1").Value = "Row ID"
$ws.Range("C1").Value = "Task"
$ws.Range("E1").Value = "Start Date"
$ws.Range("F1").Value = "End Date"

# Update the selection to reflect where the cursor ended up after the edit
$ws.Range("A2").Select()
